$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Start Build (B), Finish Build (C), and Finish Testing (G)
# for data rows 2 through 7.
for ($row = 2; $row -le 7; $row++) {
    $ws.Range("B$row").Value = 45323
    $ws.Range("C$row").Value = 45363
    $ws.Range("G$row").Value = 45364
}
